# Generate Report for Handback
# Updates the handoff/handback timestamps for the 6cf9b307-...-e2c5.md file
# (row 3 of the zh-cn / de-de report tables, row 3 of the Overview table)
# to reflect a freshly (re)generated handback report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for the 6cf9b307 file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-18 16:46:57"

# --- zh-cn sheet: Correspond Handoff / Handback datetimes for the 6cf9b307 file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-18 16:46:51"
$wsZhCn.Range("K3").Value = "2016-08-18 16:47:16"

# --- de-de sheet: Correspond Handoff / Handback datetimes for the 6cf9b307 file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-18 16:46:57"
$wsDeDe.Range("K3").Value = "2016-08-18 16:47:23"
